# The document has three inline "logo" pictures living in headers/footers:
#   - footer1.xml (default footer)      -> Pearson logo, wp:docPr id="1"
#   - footer2.xml (first-page footer)   -> Pearson logo, wp:docPr id="2"
#   - header2.xml (first-page header)   -> BTec logo,    wp:docPr id="3"
#
# The commit simply swaps each picture's display "name" with its sibling
# counterpart (image1.png <-> image2.png for the two Pearson logos, and
# image2.jpg <-> image1.jpg for the BTec logo). Word exposes this as the
# InlineShape.Name property.

$d   = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer 1 (default footer): Pearson logo -------------------------------
$shape = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
[void]$shape.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# --- Footer 2 (first-page footer): Pearson logo -----------------------------
$shape = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
[void]$shape.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# --- Header 2 (first-page header): BTec logo --------------------------------
$shape = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
[void]$shape.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"

Write-Host "Renamed the Pearson (x2) and BTec logo inline shapes"
